# Append the next quarterly evaluation row to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 84

# Column A: new date value, formatted like the existing date cells (copy style from A83)
$ws.Cells.Item($newRow, 1).Value = 45884
$ws.Range("A83").Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Column B: new value
$ws.Cells.Item($newRow, 2).Value = -0.7196185376451893
